$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 data: date serial 43852 -> 2020-01-22
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A20").Value = 43852
$ws.Range("B20").Value = "Create the question randomize for testing"
$ws.Range("D20").Value = "design the test."

# Update selection to A20 to match the saved view state
$ws.Range("A20").Select()
